# Insert a new weekly price-report row for "Poroto verde" (Vega Modelo de
# Temuco) right before the existing row 207, pushing the old rows 207-228
# down to 208-229 (matches the diff: a brand-new row 229 appears at the
# bottom, and every row from 207 on carries the data of the row above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 207, shifting rows 207:228 down to 208:229.
$ws.Rows.Item(207).EntireRow.Insert()

# Populate the newly inserted row 207 with the new record's data.
$ws.Range("A207").Value = 10
$ws.Range("B207").Value = "Vega Modelo de Temuco"
$ws.Range("C207").Value = "La Araucanía"
$ws.Range("D207").Value = 45166
$ws.Range("E207").Value = 9
$ws.Range("F207").Value = 100112031
$ws.Range("G207").Value = "Poroto verde"
$ws.Range("H207").Value = "Sin especificar"
$ws.Range("I207").Value = "Primera"
$ws.Range("J207").Value = 150
$ws.Range("K207").Value = 2000
$ws.Range("L207").Value = 2000
$ws.Range("M207").Value = 2000
$ws.Range("N207").Value = "$/kilo"
$ws.Range("O207").Value = "Provincia de Limarí"
$ws.Range("P207").Value = 2000
$ws.Range("Q207").Value = 1
$ws.Range("R207").Value = "Hortaliza"
